# NIT-9016524768.xlsx — "Elimna EC anteriores y se agregan nuevos, se modifica base de datos"
#
# The underlying worker database changed: previous account-statement (EC) rows
# are removed and replaced with a refreshed set of rows per worker (grouped by
# worker, periods 2308 -> 2302 descending), and NORVELIS DE ALBA ARRIETA's base
# salary is corrected from 1423500 to 1160000 so it matches everybody else.
#
# Layout on "Hoja1": row 15 is the header, data rows run 16..45 with
#   C = N° Doc Trabajador, D = Nombre Trabajador, E = Periodo Mora,
#   F = Valor Mora (unchanged, 46400), G = Salario Basico

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: Doc, Nombre, then the list of periods (most recent first), Salario
$workers = @(
    @{ Doc = "1148434315"; Nombre = "NORVELIS DE ALBA ARRIETA";       Periodos = @("2308","2307","2306","2305","2304");             Salario = 1160000 },
    @{ Doc = "1047455394"; Nombre = "KARINA MARGARITA MONTES CARMONA"; Periodos = @("2308","2307","2306","2305","2304","2303","2302"); Salario = 1160000 },
    @{ Doc = "1047510112"; Nombre = "MARLIZ ARRIETA JULIO";            Periodos = @("2308","2307","2306","2305","2304","2303");        Salario = 1160000 },
    @{ Doc = "1001974820"; Nombre = "NUBIS CAROLINA VERGARA SILGADO";  Periodos = @("2308","2307","2306","2305","2304","2303","2302"); Salario = 1160000 },
    @{ Doc = "1002244933"; Nombre = "ANDREA MARCELA MARIMON CORREA";   Periodos = @("2308","2306","2305","2304","2303");              Salario = 1160000 }
)

$row = 16
foreach ($worker in $workers) {
    foreach ($periodo in $worker.Periodos) {
        $ws.Range("C$row").Value = $worker.Doc
        $ws.Range("D$row").Value = $worker.Nombre
        $ws.Range("E$row").Value = $periodo
        $ws.Range("G$row").Value = $worker.Salario
        $row++
    }
}
